$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# C1 "Holdings" -> "Average"
$ws.Range("C1").Value = "Average"

# New header columns D1:F1
$ws.Range("D1").Value = "Volume"
$ws.Range("E1").Value = "Profit/Loss"
$ws.Range("F1").Value = "P/L%"

# Give the new header cells the same (bold/bordered/centered) look as the
# existing header cells by copying the format from A1 onto D1:F1.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("D1:F1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Ticker column (A2:A5) loses its bold/boxed header-style formatting ---
$ws.Range("A2:A5").Style = "Normal"

# --- Data rows ---
# Row 2 - btcinr
$ws.Range("B2").Value = 2765392
$ws.Range("C2").Value = 2199000
$ws.Range("D2").Value = 0.0001

# Row 3 - ethinr
$ws.Range("B3").Value = 120206.6
$ws.Range("C3").Value = 120000
$ws.Range("D3").Value = 0.001

# Row 4 - xrpinr
$ws.Range("B4").Value = 29.5183
$ws.Range("C4").Value = 24
$ws.Range("D4").Value = 123

# Row 5 - trxinr
$ws.Range("B5").Value = 2.4697
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 2334
